# Xbox GDK Samples - SystemInfo Readme.docx
# Update to November GDK release:
#  - remove now-stale proofing-error markers (spellStart/spellEnd, gramStart/gramEnd)
#    by collapsing the runs they used to straddle into single runs
#  - swap "Project Scarlett" wording for "an Xbox Series X|S devkit"

$d = $word.ActiveDocument

# --- Title: "SystemInfo Sample" -----------------------------------------
# Collapses the run that was wrapped in proofErr spellStart/spellEnd with
# itself, which drops the now-orphaned proofErr tags.
$d.Content.Find.Execute("SystemInfo", $true, $false, $false, $false, $false, `
    $true, 1, $false, "SystemInfo", 2) | Out-Null

# --- Description paragraph ----------------------------------------------
$d.Content.Find.Execute("This sample demonstrates a number of APIs f", $true, $false, $false, $false, $false, `
    $true, 1, $false, "This sample demonstrates a number of APIs f", 2) | Out-Null

# --- Building the sample: Xbox One devkit line ---------------------------
$d.Content.Find.Execute("If using an Xbox One devkit, set the active solution platform to Gaming.Xbox.XboxOne.x64.", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, "If using an Xbox One devkit, set the active solution platform to Gaming.Xbox.XboxOne.x64.", 2) | Out-Null

# --- Building the sample: Project Scarlett -> Xbox Series X|S devkit -----
$d.Content.Find.Execute("If using Project Scarlett, set the active solution platform to Gaming.Xbox.Scarlett.x64.", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, "If using an Xbox Series X|S devkit, set the active solution platform to Gaming.Xbox.Scarlett.x64.", 2) | Out-Null

# --- Gamepad controller instructions -------------------------------------
$d.Content.Find.Execute(", use A or DPad Right / B or DPad Left.", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", use A or DPad Right / B or DPad Left.", 2) | Out-Null

# --- Keyboard instructions -------------------------------------------------
$d.Content.Find.Execute("For keyboard, use Left or Enter / Right or BackSpace.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "For keyboard, use Left or Enter / Right or BackSpace.", 2) | Out-Null

# --- Update history: June 2020 entry --------------------------------------
$june2020 = "June 2020 " + [char]0x2013 + " Added use of GetLogicalProcessorInformation / Get LogicalProcessorInformationEx"
$d.Content.Find.Execute($june2020, $true, $false, $false, $false, $false, `
    $true, 1, $false, $june2020, 2) | Out-Null

# --- Footers: "SAMPLE: SystemInfo" / "| SystemInfo" ------------------------
foreach ($sec in $d.Sections) {
    foreach ($idx in 1, 2, 3) {
        $ftr = $sec.Footers.Item($idx)
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute("SystemInfo", $true, $false, $false, $false, $false, `
                $true, 1, $false, "SystemInfo", 2) | Out-Null
        }
    }
}
